# Remove the logo picture from the default page header (header1.xml).
# The header paragraph (styled "Kopfzeile", right-aligned) keeps its
# pPr but the run containing the <w:drawing> inline image is deleted,
# leaving an empty paragraph.

$d = $word.ActiveDocument

# Section 1's primary ("default") header corresponds to header1.xml,
# which is the one whose paragraph is right-aligned (jc="right") and
# contains the inline picture targeted by the diff.
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)  # wdHeaderFooterPrimary

if ($hdr.Exists -and $hdr.Range.InlineShapes.Count -gt 0) {
    $hdr.Range.InlineShapes.Item(1).Delete()
}
